# Weekly update: a new price record (week of 2021-10-13) was scraped for
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Apio" and inserted as
# the new row 52, pushing all subsequent rows (old 52-126) down by one
# (new 53-127). The data for every other existing row is unchanged; only
# its row number shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 52; Excel shifts rows 52:126 down to 53:127
# and carries the row-above formatting (incl. the date style on column D)
# onto the new row automatically.
$ws.Rows(52).Insert()

# Populate the newly inserted row 52 with the new weekly record.
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 44482
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 100112017
$ws.Range("G52").Value = "Apio"
$ws.Range("H52").Value = "Americana (o)"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 120
$ws.Range("K52").Value = 8000
$ws.Range("L52").Value = 9000
$ws.Range("M52").Value = 8500
$ws.Range("N52").Value = "`$/docena de matas"
$ws.Range("O52").Value = "Provincia del Elquí"
$ws.Range("P52").Value = 1417
$ws.Range("Q52").Value = 6
$ws.Range("R52").Value = "Hortaliza"
